$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-08 -> 2023-09-09, serial 45177 -> 45178) for every data row,
# from row 2 through row 372.
$lastRow = 372
$ws.Range("C2:C$lastRow").Value = 45178
